$d = $word.ActiveDocument

# --- Paragraph 1: "Project: {{.ProjectNumber}}" -- split runs w/ proofErr marks (text unchanged) ---
$d.Paragraphs.Item(1).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Project</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> {</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ProjectNumber</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p>')

# --- Paragraph 3: "Client: {{.Client}}" -- split runs w/ proofErr marks (text unchanged) ---
$d.Paragraphs.Item(3).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>Client</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>: {</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>Client</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Abadi" w:hAnsi="Abadi"/><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p>')

# --- Paragraph 4: "Status: {{.Status}}" -- split runs w/ proofErr marks (text unchanged) ---
$d.Paragraphs.Item(4).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Status: {</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Status</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p>')

# --- Paragraph 8: "Image From Tag " -> "Images From Tag " (bold) ---
$d.Paragraphs.Item(8).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>Image</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> From Tag </w:t></w:r></w:p>')

# --- Paragraph 9: "{{.Image}}" -> "PNG Image" (underline) ---
$d.Paragraphs.Item(9).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>PNG Image</w:t></w:r></w:p>')

# --- Insert two new empty paragraphs after paragraph 10 and 11 to make room for the new JPG/PNG tag blocks ---
$d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$d.Paragraphs.Item(11).Range.InsertParagraphAfter()

# --- Paragraph 10 (was empty): add "{{.ImagePng}}" tag with proofErr marks ---
$d.Paragraphs.Item(10).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Image</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Png</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p>')

# --- Paragraph 11 (new): "JPG Image" (underline) ---
$d.Paragraphs.Item(11).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>JPG</w:t></w:r><w:r><w:rPr><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Image</w:t></w:r></w:p>')

# --- Paragraph 12 (new): add "{{.ImageJpg}}" tag with proofErr marks ---
$d.Paragraphs.Item(12).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Image</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Jp</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>g</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>}}</w:t></w:r></w:p>')

# --- Paragraph 13 (was "End of doc`"): now an empty paragraph ---
$d.Paragraphs.Item(13).Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>')

